$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# GeneNames (column E) reorderings
$ws.Range("E14").Value = "PCDHGA10;PCDHGA9;PCDHGA4;PCDHGA3;PCDHGA8;PCDHGB5;PCDHGA12;PCDHGA2;PCDHGB6;PCDHGB4;PCDHGB7;PCDHGB3;PCDHGA11;PCDHGA5;PCDHGB1;PCDHGB2;PCDHGA7;PCDHGA6;PCDHGA1"
$ws.Range("E19").Value = "RP11-550P17.5;RP11-180D21.3"
$ws.Range("E35").Value = "PMM2;RP11-152P23.2;RP11-77H9.2"
$ws.Range("E48").Value = "GSN-AS1;GSN"
$ws.Range("E51").Value = "C16orf45;RP11-1021N1.1"
$ws.Range("E70").Value = "PMM2;RP11-152P23.2;RP11-77H9.2"
$ws.Range("E76").Value = "MYADM;AC008753.6"
$ws.Range("E80").Value = "H1FX-AS1;H1FX"
$ws.Range("E98").Value = "FGF11;RP11-104H15.10;RP11-104H15.7;RP11-104H15.8"
$ws.Range("E101").Value = "GPR1;GPR1-AS"
$ws.Range("E113").Value = "BZRAP1-AS1;RNF43"

# GeneClasses (column F) reorderings
$ws.Range("F19").Value = "lincRNA;antisense"
$ws.Range("F35").Value = "protein_coding;antisense"
$ws.Range("F48").Value = "protein_coding;antisense"
$ws.Range("F54").Value = "protein_coding;miRNA"
$ws.Range("F70").Value = "protein_coding;antisense"
$ws.Range("F76").Value = "protein_coding;antisense"
$ws.Range("F80").Value = "protein_coding;antisense"
$ws.Range("F98").Value = "protein_coding;antisense;processed_transcript"
$ws.Range("F101").Value = "protein_coding;processed_transcript"
$ws.Range("F113").Value = "protein_coding;antisense"
